# Auto-generated cell updates applying the crypto price refresh diff
# (GitHub Actions scheduled refresh of coinranking.com snapshot data).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, $value) {
    # Force text storage so numeric-looking strings (e.g. "1.003",
    # "30.110.65") are not silently reinterpreted as numbers/dates,
    # matching the inlineStr cells in the source workbook.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = '30.110.65'
$ws.Range("E2").Value = '  +0.92%  '

$ws.Range("D3").Value = '1.890.14'
$ws.Range("E3").Value = '  +0.13%  '

Set-TextCell $ws.Range("D4") '1.003'
$ws.Range("E4").Value = '  +0.26%  '

Set-TextCell $ws.Range("D5") '0.7384'
$ws.Range("E5").Value = '  -1.98%  '

Set-TextCell $ws.Range("D6") '242.76'
$ws.Range("E6").Value = '  +0.16%  '

Set-TextCell $ws.Range("D7") '1.002'
$ws.Range("E7").Value = '  +0.18%  '

Set-TextCell $ws.Range("D8") '0.3175'

Set-TextCell $ws.Range("D9") '24.84'
$ws.Range("E9").Value = '  -1.96%  '

Set-TextCell $ws.Range("D10") '0.07147'
$ws.Range("E10").Value = '  +0.41%  '

Set-TextCell $ws.Range("D11") '0.08341'
$ws.Range("E11").Value = '  -1.61%  '

$ws.Range("D12").Value = '1.946.61'
$ws.Range("E12").Value = '  +3.35%  '

Set-TextCell $ws.Range("D13") '0.7580'
$ws.Range("E13").Value = '  -0.25%  '

Set-TextCell $ws.Range("D14") '5.418'
$ws.Range("E14").Value = '  +1.04%  '

Set-TextCell $ws.Range("D15") '93.05'
$ws.Range("E15").Value = '  -0.35%  '

Set-TextCell $ws.Range("D16") '6.159'
$ws.Range("E16").Value = '  +0.34%  '

$ws.Range("D17").Value = '30.102.79'
$ws.Range("E17").Value = '  +1.15%  '

Set-TextCell $ws.Range("D18") '251.17'
$ws.Range("E18").Value = '  +3.04%  '

Set-TextCell $ws.Range("D19") '13.60'
$ws.Range("E19").Value = '  -0.79%  '

Set-TextCell $ws.Range("D20") '0.000007860'
$ws.Range("E20").Value = '  +0.68%  '

$ws.Range("D21").Value = '2.198.18'
$ws.Range("E21").Value = '  +2.85%  '

$ws.Range("E22").Value = '  +0.14%  '

Set-TextCell $ws.Range("D23") '7.930'
$ws.Range("E23").Value = '  -0.75%  '

Set-TextCell $ws.Range("D24") '1.002'
$ws.Range("E24").Value = '  +0.14%  '

Set-TextCell $ws.Range("D25") '0.1571'
$ws.Range("E25").Value = '  -1.65%  '

Set-TextCell $ws.Range("D26") '9.306'
$ws.Range("E26").Value = '  -0.76%  '

Set-TextCell $ws.Range("D27") '164.69'
$ws.Range("E27").Value = '  +1.27%  '

Set-TextCell $ws.Range("D28") '18.72'
$ws.Range("E28").Value = '  +0.01%  '

Set-TextCell $ws.Range("D29") '2.057'
$ws.Range("E29").Value = '  +1.28%  '

Set-TextCell $ws.Range("D30") '1.480'
$ws.Range("E30").Value = '  -0.49%  '

Set-TextCell $ws.Range("D31") '4.577'
$ws.Range("E31").Value = '  +1.38%  '

Set-TextCell $ws.Range("D32") '1.540'
$ws.Range("E32").Value = '  +0.17%  '

Set-TextCell $ws.Range("D33") '4.197'
$ws.Range("E33").Value = '  +1.65%  '

Set-TextCell $ws.Range("D34") '0.05351'
$ws.Range("E34").Value = '  -1.28%  '

Set-TextCell $ws.Range("D35") '1.255'
$ws.Range("E35").Value = '  +1.00%  '

Set-TextCell $ws.Range("D36") '0.7706'
$ws.Range("E36").Value = '  +2.69%  '

Set-TextCell $ws.Range("D37") '1.003'
$ws.Range("E37").Value = '  +0.15%  '

Set-TextCell $ws.Range("D38") '2.728'
$ws.Range("E38").Value = '  +0.67%  '

$ws.Range("E39").Value = '  +0.75%  '

$ws.Range("E40").Value = '  -0.38%  '

Set-TextCell $ws.Range("D41") '0.4563'
$ws.Range("E41").Value = '  +2.37%  '

$ws.Range("D42").Value = '1.101.82'
$ws.Range("E42").Value = '  +0.96%  '

Set-TextCell $ws.Range("D43") '6.071'
$ws.Range("E43").Value = '  -0.54%  '

Set-TextCell $ws.Range("D44") '72.55'
$ws.Range("E44").Value = '  +0.01%  '

Set-TextCell $ws.Range("D45") '0.8750'
$ws.Range("E45").Value = '  +1.68%  '

Set-TextCell $ws.Range("D46") '104.32'
$ws.Range("E46").Value = '  +1.94%  '

Set-TextCell $ws.Range("D47") '1.003'
$ws.Range("E47").Value = '  +0.26%  '

Set-TextCell $ws.Range("D48") '1.862'
$ws.Range("E48").Value = '  +0.18%  '

Set-TextCell $ws.Range("D49") '7.583'
$ws.Range("E49").Value = '  -1.85%  '

$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell $ws.Range("D50") '9.613'
$ws.Range("E50").Value = '  -1.49%  '

$ws.Range("B51").Value = 'RocketPoolETH'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D51").Value = '2.062.95'
$ws.Range("E51").Value = '  +1.45%  '
